$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style) from A16 into the new rows A17:A19
$ws.Range("A16").Copy()
$ws.Range("A17:A19").PasteSpecial(-4122)

# Row 10: Gaussian-Quadrature
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("C10").Value = 1.07982118082457
$ws.Range("D10").Value = 0.7630299429227828
$ws.Range("E10").Value = 1.035803132436836
$ws.Range("F10").Value = 1.07982118082457
$ws.Range("G10").Value = 0.8744127469823707
$ws.Range("H10").Value = 1.100669706456135
$ws.Range("I10").Value = 1.050204621539189
$ws.Range("J10").Value = 0.7630299429227828
$ws.Range("K10").Value = 0.8994165376798096
$ws.Range("L10").Value = 0.9896188592521895
$ws.Range("M10").Value = 0.9839902218603139

# Row 11: Spiral-90deg-10rot-5space
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value = 0.9986316039002782
$ws.Range("D11").Value = 1.014515462032718
$ws.Range("E11").Value = 0.9815533929044289
$ws.Range("F11").Value = 0.9986316039002782
$ws.Range("G11").Value = 1.017007697569675
$ws.Range("H11").Value = 0.9393484151487657
$ws.Range("I11").Value = 0.9852205225525594
$ws.Range("J11").Value = 1.014515462032718
$ws.Range("K11").Value = 0.9980344274685736
$ws.Range("L11").Value = 0.9983330156844259
$ws.Range("M11").Value = 0.9893795156847376

# Row 12: Spiral-90deg-15rot-5space
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value = 0.9977409702172951
$ws.Range("D12").Value = 1.015679590969903
$ws.Range("E12").Value = 0.9815666993021648
$ws.Range("F12").Value = 0.9977409702172951
$ws.Range("G12").Value = 1.017596285995758
$ws.Range("H12").Value = 0.9392084747018541
$ws.Range("I12").Value = 0.9849176187269965
$ws.Range("J12").Value = 1.015679590969903
$ws.Range("K12").Value = 0.9986231451360339
$ws.Range("L12").Value = 0.9981820576766645
$ws.Range("M12").Value = 0.9894516066523286

# Row 13: Spiral-90deg-10rot-3space
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value = 0.9985434333057849
$ws.Range("D13").Value = 1.01461195381169
$ws.Range("E13").Value = 0.9815906552222895
$ws.Range("F13").Value = 0.9985434333057849
$ws.Range("G13").Value = 1.017092012722028
$ws.Range("H13").Value = 0.939252658025519
$ws.Range("I13").Value = 0.985133181604265
$ws.Range("J13").Value = 1.01461195381169
$ws.Range("K13").Value = 0.9981013045169895
$ws.Range("L13").Value = 0.9983223689113871
$ws.Range("M13").Value = 0.9893706491152626

# Row 14: NoRotation-tilt60deg
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("C14").Value = 0.9581080000000002
$ws.Range("D14").Value = 0.6402800000000011
$ws.Range("E14").Value = 1.1411
$ws.Range("F14").Value = 0.9581080000000002
$ws.Range("G14").Value = 0.7521880000000006
$ws.Range("H14").Value = 1.451455999999997
$ws.Range("I14").Value = 1.103799999999996
$ws.Range("J14").Value = 0.6402800000000011
$ws.Range("K14").Value = 0.8906900000000008
$ws.Range("L14").Value = 0.9243990000000004
$ws.Range("M14").Value = 1.007821999999999

# Row 15: Rotation-NoTilt
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("C15").Value = 0.98
$ws.Range("D15").Value = 0.26
$ws.Range("E15").Value = 1.27
$ws.Range("F15").Value = 0.98
$ws.Range("G15").Value = 0.5
$ws.Range("H15").Value = 1.87
$ws.Range("I15").Value = 1.22
$ws.Range("J15").Value = 0.26
$ws.Range("K15").Value = 0.765
$ws.Range("L15").Value = 0.8724999999999999
$ws.Range("M15").Value = 1.016666666666667

# Row 16: Rotation-60detTilt
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("C16").Value = 0.9893526728704014
$ws.Range("D16").Value = 0.5656516016128018
$ws.Range("E16").Value = 1.155065619251198
$ws.Range("F16").Value = 0.9893526728704014
$ws.Range("G16").Value = 0.7081621151744009
$ws.Range("H16").Value = 1.497412810342399
$ws.Range("I16").Value = 1.123348275404795
$ws.Range("J16").Value = 0.5656516016128018
$ws.Range("K16").Value = 0.8603586104319998
$ws.Range("L16").Value = 0.9248556416512006
$ws.Range("M16").Value = 1.006498849109333

# Row 17: HexGrid-90degTilt5degRes
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value = 0.9939584094867431
$ws.Range("D17").Value = 0.994647621757706
$ws.Range("E17").Value = 0.9950709390557989
$ws.Range("F17").Value = 0.9939584094867431
$ws.Range("G17").Value = 0.9929422555134165
$ws.Range("H17").Value = 0.9954400091541379
$ws.Range("I17").Value = 0.9947444320283918
$ws.Range("J17").Value = 0.994647621757706
$ws.Range("K17").Value = 0.9948592804067524
$ws.Range("L17").Value = 0.9944088449467479
$ws.Range("M17").Value = 0.994467277832699

# Row 18: HexGrid-90degTilt22p5degRes
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value = 0.9678035059418442
$ws.Range("D18").Value = 1.033984559192438
$ws.Range("E18").Value = 0.9949299499943953
$ws.Range("F18").Value = 0.9678035059418442
$ws.Range("G18").Value = 1.010033665185126
$ws.Range("H18").Value = 0.9938803833157001
$ws.Range("I18").Value = 0.9873534430311675
$ws.Range("J18").Value = 1.033984559192438
$ws.Range("K18").Value = 1.014457254593417
$ws.Range("L18").Value = 0.9911303802676303
$ws.Range("M18").Value = 0.9979975844434451

# Row 19: HexGrid-60degTilt5degRes
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value = 0.9762665294871865
$ws.Range("D19").Value = 1.090996924008518
$ws.Range("E19").Value = 0.9701505717504328
$ws.Range("F19").Value = 0.9762665294871865
$ws.Range("G19").Value = 1.04952550139501
$ws.Range("H19").Value = 0.9237341146239977
$ws.Range("I19").Value = 0.967981570227083
$ws.Range("J19").Value = 1.090996924008518
$ws.Range("K19").Value = 1.030573747879476
$ws.Range("L19").Value = 1.003420138683331
$ws.Range("M19").Value = 0.9964425352487049
